# Auto-generated edit script applying scheduled-runner value updates
# to the Belias_Profits workbook (per-sheet, per-cell numeric updates).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1275.7
$ws.Range("I19").Value = 827.6667
$ws.Range("J19").Value = 1947.75
$ws.Range("K19").Value = 827.6667
$ws.Range("L19").Value = 1947.75
$ws.Range("M19").Value = -652.6667
$ws.Range("N19").Value = -2297.75
$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -31996
$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -99984
$ws.Range("H111").Value = 5765.5
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 5765.5
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 17296.5
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -23430.5
$ws.Range("H113").Value = 4256.2856
$ws.Range("I113").Value = 4179.6523
$ws.Range("J113").Value = 4403.1665
$ws.Range("K113").Value = 4179.6523
$ws.Range("L113").Value = 4403.1665
$ws.Range("M113").Value = -925.6522999999997
$ws.Range("N113").Value = -10911.1665
$ws.Range("H127").Value = 333333730
$ws.Range("I127").Value = 333333730
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 1000001190
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -999996230
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 15154015
$ws.Range("I45").Value = 15154015
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 15154015
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -15153638
$ws.Range("N45").ClearContents()
$ws.Range("H46").Value = 127238
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 127238
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 127238
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -127876
$ws.Range("H52").Value = 13199.8
$ws.Range("J52").Value = 13199.8
$ws.Range("L52").Value = 13199.8
$ws.Range("N52").Value = -13835.8
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30676
$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32340

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2267.8
$ws.Range("I94").Value = 1628.8182
$ws.Range("J94").Value = 4025
$ws.Range("K94").Value = 1628.8182
$ws.Range("L94").Value = 4025
$ws.Range("M94").Value = -1177.8182
$ws.Range("N94").Value = -4927
$ws.Range("H107").Value = 1431.1666
$ws.Range("I107").Value = 1286.1
$ws.Range("J107").Value = 1612.5
$ws.Range("K107").Value = 1286.1
$ws.Range("L107").Value = 1612.5
$ws.Range("M107").Value = 633.9000000000001
$ws.Range("N107").Value = -5452.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2121.4688
$ws.Range("I31").Value = 1240.2778
$ws.Range("J31").Value = 3254.4285
$ws.Range("K31").Value = 1240.2778
$ws.Range("L31").Value = 3254.4285
$ws.Range("M31").Value = -945.2778000000001
$ws.Range("N31").Value = -3844.4285
$ws.Range("H34").Value = 2121.4688
$ws.Range("I34").Value = 1240.2778
$ws.Range("J34").Value = 3254.4285
$ws.Range("K34").Value = 1240.2778
$ws.Range("L34").Value = 3254.4285
$ws.Range("M34").Value = -1038.2778
$ws.Range("N34").Value = -3658.4285
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H99").Value = 4877.75
$ws.Range("I99").Value = 4877.75
$ws.Range("K99").Value = 4877.75
$ws.Range("M99").Value = -3379.75
$ws.Range("H126").Value = 4877.75
$ws.Range("I126").Value = 4877.75
$ws.Range("K126").Value = 14633.25
$ws.Range("M126").Value = -12163.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1721.4286
$ws.Range("I17").Value = 666.44446
$ws.Range("J17").Value = 3620.4
$ws.Range("K17").Value = 1999.33338
$ws.Range("L17").Value = 10861.2
$ws.Range("M17").Value = -1830.33338
$ws.Range("N17").Value = -11199.2
$ws.Range("H34").Value = 1045.6923
$ws.Range("I34").Value = 546
$ws.Range("J34").Value = 1267.7778
$ws.Range("K34").Value = 1638
$ws.Range("L34").Value = 3803.3334
$ws.Range("M34").Value = -1554
$ws.Range("N34").Value = -3971.3334
$ws.Range("H39").Value = 2627.5881
$ws.Range("I39").Value = 1033.3334
$ws.Range("J39").Value = 2969.2144
$ws.Range("K39").Value = 3100.0002
$ws.Range("L39").Value = 8907.643199999999
$ws.Range("M39").Value = -2806.0002
$ws.Range("N39").Value = -9495.643199999999
$ws.Range("H55").Value = 795.25
$ws.Range("I55").Value = 441
$ws.Range("J55").Value = 913.3333
$ws.Range("K55").Value = 1323
$ws.Range("L55").Value = 2739.9999
$ws.Range("M55").Value = -1146
$ws.Range("N55").Value = -3093.9999
$ws.Range("H75").Value = 4383.25
$ws.Range("I75").Value = 1903.25
$ws.Range("J75").Value = 5209.9165
$ws.Range("K75").Value = 5709.75
$ws.Range("L75").Value = 15629.7495
$ws.Range("M75").Value = -4711.75
$ws.Range("N75").Value = -17625.7495
$ws.Range("H78").Value = 4383.25
$ws.Range("I78").Value = 1903.25
$ws.Range("J78").Value = 5209.9165
$ws.Range("K78").Value = 17129.25
$ws.Range("L78").Value = 46889.2485
$ws.Range("M78").Value = -12137.25
$ws.Range("N78").Value = -56873.2485
$ws.Range("H131").Value = 942.8659699999999
$ws.Range("J131").Value = 943.3125
$ws.Range("L131").Value = 2829.9375
$ws.Range("N131").Value = -12909.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H103").Value = 24300
$ws.Range("J103").Value = 24300
$ws.Range("L103").Value = 24300
$ws.Range("N103").Value = -26644

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 9750
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 9750
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 9750
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -10564
$ws.Range("H120").Value = 30697.9
$ws.Range("J120").Value = 30697.9
$ws.Range("L120").Value = 30697.9
$ws.Range("N120").Value = -40373.9
$ws.Range("H132").Value = 3579.6726
$ws.Range("I132").Value = 3552.7942
$ws.Range("J132").Value = 3623.1904
$ws.Range("K132").Value = 10658.3826
$ws.Range("L132").Value = 10869.5712
$ws.Range("M132").Value = -8128.382599999999
$ws.Range("N132").Value = -15929.5712

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 7192.3335
$ws.Range("I107").Value = 10378.5
$ws.Range("J107").Value = 820
$ws.Range("K107").Value = 31135.5
$ws.Range("L107").Value = 2460
$ws.Range("M107").Value = -29215.5
$ws.Range("N107").Value = -6300
$ws.Range("H120").Value = 29051.25
$ws.Range("J120").Value = 29051.25
$ws.Range("L120").Value = 29051.25
$ws.Range("N120").Value = -38727.25

Write-Host "Applied 184 cell updates and 9 cell removals across 8 sheets."
